$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'26.094.56"
$ws.Range("E2").Value = "  -4.47%  "
$ws.Range("D3").Value = "'1.649.82"
$ws.Range("E3").Value = "  -3.54%  "
$ws.Range("E4").Value = "  -0.21%  "
$ws.Range("D5").Value = "'215.66"
$ws.Range("E5").Value = "  -3.86%  "
$ws.Range("D6").Value = "'0.5106"
$ws.Range("E6").Value = "  -2.97%  "
$ws.Range("E7").Value = "  -0.04%  "
$ws.Range("D8").Value = "'0.2592"
$ws.Range("E8").Value = "  -1.59%  "
$ws.Range("E9").Value = "  -2.82%  "
$ws.Range("D10").Value = "'19.67"
$ws.Range("E10").Value = "  -4.73%  "
$ws.Range("E11").Value = "  -0.35%  "
$ws.Range("D12").Value = "'1.645.06"
$ws.Range("E12").Value = "  -3.63%  "
$ws.Range("D13").Value = "'4.276"
$ws.Range("E13").Value = "  -3.62%  "
$ws.Range("D14").Value = "'1.878.60"
$ws.Range("E14").Value = "  -3.53%  "
$ws.Range("D15").Value = "'0.5476"
$ws.Range("E15").Value = "  -4.98%  "
$ws.Range("D16").Value = "'0.0₅7996"
$ws.Range("D17").Value = "'63.88"
$ws.Range("E17").Value = "  -5.30%  "
$ws.Range("D18").Value = "'26.111.87"
$ws.Range("E18").Value = "  -4.55%  "
$ws.Range("E19").Value = "  +0.01%  "
$ws.Range("D20").Value = "'207.55"
$ws.Range("E20").Value = "  -4.85%  "
$ws.Range("D21").Value = "'4.384"
$ws.Range("E21").Value = "  -5.31%  "
$ws.Range("E22").Value = "  -3.35%  "
$ws.Range("D23").Value = "'6.012"
$ws.Range("E23").Value = "  -0.08%  "
$ws.Range("D24").Value = "'1.009"
$ws.Range("E24").Value = "  -0.01%  "
$ws.Range("D25").Value = "'1.873"
$ws.Range("E25").Value = "  +8.54%  "
$ws.Range("D26").Value = "'143.14"
$ws.Range("E26").Value = "  -1.47%  "
$ws.Range("E27").Value = "  -2.41%  "
$ws.Range("D28").Value = "'6.916"
$ws.Range("E28").Value = "  -3.76%  "
$ws.Range("D29").Value = "'15.82"
$ws.Range("E29").Value = "  -1.79%  "
$ws.Range("D30").Value = "'0.05067"
$ws.Range("E30").Value = "  -4.31%  "
$ws.Range("D31").Value = "'1.244"
$ws.Range("E31").Value = "  -3.78%  "
$ws.Range("D32").Value = "'3.343"
$ws.Range("E32").Value = "  -3.33%  "
$ws.Range("D33").Value = "'3.232"
$ws.Range("E33").Value = "  -3.29%  "
$ws.Range("D34").Value = "'1.549"
$ws.Range("E34").Value = "  -5.27%  "
$ws.Range("D35").Value = "'2.348"
$ws.Range("E35").Value = "  -2.30%  "
$ws.Range("D36").Value = "'0.9148"
$ws.Range("E36").Value = "  -3.39%  "
$ws.Range("D37").Value = "'2.649"
$ws.Range("E37").Value = "  -6.45%  "
$ws.Range("D38").Value = "'0.5702"
$ws.Range("E38").Value = "  -2.71%  "
$ws.Range("D39").Value = "'1.144.80"
$ws.Range("E39").Value = "  -2.87%  "
$ws.Range("D40").Value = "'0.01575"
$ws.Range("E40").Value = "  -4.30%  "
$ws.Range("D41").Value = "'2.563"
$ws.Range("E41").Value = "  -0.22%  "
$ws.Range("E42").Value = "  -0.03%  "
$ws.Range("D43").Value = "'5.663"
$ws.Range("E43").Value = "  -1.75%  "
$ws.Range("D44").Value = "'0.8247"
$ws.Range("E44").Value = "  -1.60%  "
$ws.Range("D45").Value = "'100.08"
$ws.Range("E45").Value = "  -0.87%  "
$ws.Range("D46").Value = "'1.789.92"
$ws.Range("E46").Value = "  -3.55%  "
$ws.Range("D47").Value = "'0.0₈114"
$ws.Range("E47").Value = "  -3.47%  "
$ws.Range("D48").Value = "'0.4537"
$ws.Range("E48").Value = "  -0.56%  "
$ws.Range("D49").Value = "'1.011"
$ws.Range("E49").Value = "  +0.54%  "
$ws.Range("D50").Value = "'55.25"
$ws.Range("E50").Value = "  -3.69%  "
$ws.Range("D51").Value = "'7.811"
$ws.Range("E51").Value = "  -3.84%  "
